$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 4 for rule R3 (rule id first) ------------------------------
$ws.Cells.Item(4, 1).Value2 = "R3"

# --- Update existing "Solution" texts for R1 (row2) and R2 (row3) -----------
$ws.Cells.Item(2, 3).Value2 = "Solution1: R1:`n1. Ask client for microsoft error debugging."
$ws.Cells.Item(3, 3).Value2 = "Solution2: R2: `n1. Ask client for Security token ."

# --- Fill in the rest of the new row 4 --------------------------------------
$ws.Cells.Item(4, 2).Value2 = "[{`n`t""message"": ""VALUEADDCO""`n}]"
$ws.Cells.Item(4, 3).Value2 = "Solution2: R3:`n1. Problem is in VALUEADDCO, Ask client to change the value."

# Copy the formatting (wrap text style) from row 3 to the new row 4 cells
$ws.Range("B3:C3").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(4).RowHeight = 95.25

# --- Add an extra (mostly empty) row 5, matching formatting of column B -----
$ws.Cells.Item(5, 2).Value2 = $null
$ws.Cells.Item(5, 2).WrapText = $true

# --- Update the view so the new rows are visible ----------------------------
$excel.ActiveWindow.ScrollRow = 4
$ws.Cells.Item(4, 2).Select()
